# "Generate Report for Handback" - refresh the localization-status report
# after a successful handback (target is now in sync with en-US).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status cells for the localized file ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# --- zh-cn sheet ---
# Status column (C) for both tracked files now reads "Handed back..."
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

# Latest Handback DateTime (K) refreshed to the new handback run timestamp
$wsZhCn.Range("K2").Value = "2016-10-19 17:57:01"
$wsZhCn.Range("K3").Value = "2016-10-19 17:57:01"

# Error Detail (P) no longer applicable - handback file is now current
$wsZhCn.Range("P2").Value = ""

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

$wsDeDe.Range("K2").Value = "2016-10-19 17:57:19"
$wsDeDe.Range("K3").Value = "2016-10-19 17:57:19"

$wsDeDe.Range("P2").Value = ""

# --- Column widths: the longer status text / cleared error text change the
# natural auto-fit width of the affected columns ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333334

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333334
